# Quarterly indexing bug-fix: the qoq-error series for each origin-quarter
# row was being written one column short (starting at the 2nd horizon
# instead of the 1st). Fix: shift each row's existing error values one
# column to the right (B->C, C->D, ... ) and write the correct first-horizon
# value into column B. The last existing trailing value in rows that already
# reached column K falls off (the sheet only has columns A:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-B (first-horizon) values to insert for each data row (2-16).
$newFirstValues = @{
    2  = -0.5825945370336409
    3  = 0.09567504080935779
    4  = -0.2604190369987228
    5  = 0.8354549961584912
    6  = -0.1000793599026215
    7  = -0.3537865060796963
    8  = 0.1481773904324453
    9  = 0.157445989004155
    10 = -0.5006594565260708
    11 = 0.2803578805354692
    12 = -0.1719748578450117
    13 = 0.3058625397463315
    14 = -0.6123299526872862
    15 = 0.6883713851991116
    16 = -0.2766911554241067
}

$lastDataCol = 11   # column K

for ($r = 2; $r -le 16; $r++) {
    # Find how many contiguous data cells currently exist starting at column B.
    $lastCol = 1   # column A
    for ($c = 2; $c -le $lastDataCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -ne $null) {
            $lastCol = $c
        }
    }

    if ($lastCol -ge 2) {
        # Shift existing values one column to the right, right-to-left so we
        # don't clobber a value before it's been read. Anything that would
        # land past column K is simply dropped (matches the fixed sheet width).
        $destLastCol = [Math]::Min($lastCol + 1, $lastDataCol)
        for ($c = $destLastCol; $c -ge 3; $c--) {
            $srcVal = $ws.Cells.Item($r, $c - 1).Value2
            $ws.Cells.Item($r, $c).Value = $srcVal
        }
    }

    # Write the new first-horizon value into column B.
    $ws.Cells.Item($r, 2).Value = $newFirstValues[$r]
}
